$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalidCredentialTest")

# Additional "peter" login rows (invalid-credential test data), appended
# below the existing peter123 row - modelled as an object array of rows
# in the calling test, each row written out cell-by-cell here.
$newRows = @(
    @("peter", "peter124", "Dutch", "Invalid username or password"),
    @("peter", "peter125", "Dutch", "Invalid username or password"),
    @("peter", "peter126", "Dutch", "Invalid username or password"),
    @("peter", "peter127", "Dutch", "Invalid username or password"),
    @("peter", "peter128", "Dutch", "Invalid username or password")
)

$r = 4
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Resize columns to fit their (now longer) contents.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Restore the selection to A3, as left by the author.
$ws.Range("A3").Select()
